$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    4368.340131353961,
    4352.320777739636,
    4343.541579856456,
    4288.947294527362,
    4288.947294527362,
    4260.964087326761,
    4193.301722847783,
    4193.301722847783,
    4193.301722847783,
    4036.53995309705,
    4036.53995309705
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $values[$i]
}
